# Apply updated cryptocurrency price/volume data to the worksheet.
# Numeric-looking "Price" values are prefixed with a leading apostrophe so
# Excel keeps them as text (matching the workbook's original string data)
# instead of auto-converting them into Number cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.348.75'
$ws.Range('E2').Value = '  +5.67%  '
$ws.Range('D3').Value = '2.743.92'
$ws.Range('E3').Value = '  +3.50%  '
$ws.Range('D4').Value = '''0.998'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '''581.64'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').Value = '''157.95'
$ws.Range('E6').Value = '  +9.30%  '
$ws.Range('D7').Value = '''0.996'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = '''0.611'
$ws.Range('E8').Value = '  +1.99%  '
$ws.Range('D9').Value = '2.772.78'
$ws.Range('E9').Value = '  +3.96%  '
$ws.Range('D10').Value = '''6.78'
$ws.Range('E10').Value = '  +3.31%  '
$ws.Range('E11').Value = '  +5.88%  '
$ws.Range('D12').Value = '''0.393'
$ws.Range('E12').Value = '  +2.90%  '
$ws.Range('E13').Value = '  +2.54%  '
$ws.Range('D14').Value = '3.235.14'
$ws.Range('E14').Value = '  +3.57%  '
$ws.Range('D15').Value = '''26.94'
$ws.Range('E15').Value = '  +3.08%  '
$ws.Range('D16').Value = '64.206.88'
$ws.Range('E16').Value = '  +5.50%  '
$ws.Range('E17').Value = '  +7.35%  '
$ws.Range('D18').Value = '2.766.08'
$ws.Range('E18').Value = '  +3.97%  '
$ws.Range('D19').Value = '''12.07'
$ws.Range('E19').Value = '  +3.93%  '
$ws.Range('D20').Value = '''4.94'
$ws.Range('E20').Value = '  +4.48%  '
$ws.Range('D21').Value = '''363.52'
$ws.Range('E21').Value = '  +3.76%  '
$ws.Range('D22').Value = '''7.03'
$ws.Range('E22').Value = '  +1.35%  '
$ws.Range('D23').Value = '''0.997'
$ws.Range('E23').Value = '  -0.36%  '
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('D25').Value = '''66.90'
$ws.Range('E25').Value = '  +4.59%  '
$ws.Range('E26').Value = '  +5.80%  '
$ws.Range('D27').Value = '''8.61'
$ws.Range('E27').Value = '  +5.53%  '
$ws.Range('D28').Value = '''0.999'
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('D29').Value = '0.0₃0911'
$ws.Range('E29').Value = '  +12.25%  '
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('E31').Value = '  +4.81%  '
$ws.Range('D32').Value = '''1.28'
$ws.Range('E32').Value = '  +20.50%  '
$ws.Range('D33').Value = '''174.11'
$ws.Range('E33').Value = '  +4.76%  '
$ws.Range('D34').Value = '''0.996'
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('D35').Value = '''20.67'
$ws.Range('E35').Value = '  +3.89%  '
$ws.Range('D36').Value = '''4.90'
$ws.Range('E36').Value = '  +7.34%  '
$ws.Range('E37').Value = '  +8.72%  '
$ws.Range('D38').Value = '''1.83'
$ws.Range('E38').Value = '  +10.50%  '
$ws.Range('E39').Value = '  +12.69%  '
$ws.Range('D40').Value = '''343.64'
$ws.Range('E40').Value = '  +1.01%  '
$ws.Range('D41').Value = '''4.29'
$ws.Range('E41').Value = '  +5.65%  '
$ws.Range('D42').Value = '''39.33'
$ws.Range('E42').Value = '  +1.94%  '
$ws.Range('D43').Value = '''5.89'
$ws.Range('E43').Value = '  +11.89%  '
$ws.Range('D44').Value = '''22.12'
$ws.Range('E44').Value = '  +8.44%  '
$ws.Range('D45').Value = '''22.13'
$ws.Range('E45').Value = '  +7.34%  '
$ws.Range('D46').Value = '''0.0597'
$ws.Range('E46').Value = '  +6.38%  '
$ws.Range('D47').Value = '''0.652'
$ws.Range('E47').Value = '  +5.36%  '
$ws.Range('D48').Value = '''138.33'
$ws.Range('E48').Value = '  +3.26%  '
$ws.Range('E49').Value = '  +3.85%  '
$ws.Range('E50').Value = '  +2.25%  '
$ws.Range('E51').Value = '  -0.01%  '
